$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-11 (row 24)
$ws.Range("B24").Value = 6385
$ws.Range("D24").Value = 5965698
$ws.Range("E24").Value = 934.3301487862177
$ws.Range("F24").Value = 8.847596317763372
$ws.Range("H24").Value = 26.3765916503339
